# Generate Report for Handoff
# Updates status text and timestamps on each sheet, and narrows the
# "datetime" columns on the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Update timestamps (stored as plain text, not as real dates)
$wsOverview.Range("G2").Value = "2016-08-26 12:59:41"
$wsDeDe.Range("H2").Value = "2016-08-26 12:59:41"
$wsZhCn.Range("H2").Value = "2016-08-26 12:59:37"

# --- Narrow the datetime columns (target raw width 17.2159881591797 sits
#     between Excel's quantized column-width grid points; 16.335 is the
#     ColumnWidth value that lands on the closest achievable grid point)
$wsOverview.Columns.Item(5).ColumnWidth = 16.335
$wsOverview.Columns.Item(6).ColumnWidth = 16.335
$wsZhCn.Columns.Item(3).ColumnWidth = 16.335
$wsDeDe.Columns.Item(3).ColumnWidth = 16.335
